$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.932.29"
$ws.Range("E2").Value = "  -1.25%  "

$ws.Range("D3").Value = "1.911.28"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'324.28"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").Value = "'0.4593"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").Value = "'0.3822"
$ws.Range("E8").Value = "  -1.34%  "

$ws.Range("D9").Value = "'0.07704"
$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("D10").Value = "'0.9805"
$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("D11").Value = "'22.22"
$ws.Range("E11").Value = "  -2.03%  "

$ws.Range("D12").Value = "1.865.48"
$ws.Range("E12").Value = "  -4.55%  "

$ws.Range("D13").Value = "'5.688"
$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("D14").Value = "'6.955"
$ws.Range("E14").Value = "  -1.57%  "

$ws.Range("D15").Value = "'0.07056"
$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "'83.86"
$ws.Range("E17").Value = "  -3.31%  "

$ws.Range("D18").Value = "'0.000009452"
$ws.Range("E18").Value = "  -2.87%  "

$ws.Range("D19").Value = "'16.62"
$ws.Range("E19").Value = "  -2.23%  "

$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").Value = "28.923.89"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").Value = "'5.316"
$ws.Range("E22").Value = "  -2.88%  "

$ws.Range("E23").Value = "  -1.48%  "

$ws.Range("D24").Value = "'2.098"
$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").Value = "'158.50"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").Value = "'19.02"
$ws.Range("E26").Value = "  -1.87%  "

$ws.Range("D27").Value = "'5.683"
$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("D28").Value = "'117.74"
$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("D29").Value = "'1.878"
$ws.Range("E29").Value = "  +2.47%  "

$ws.Range("D30").Value = "'0.09286"
$ws.Range("E30").Value = "  -0.67%  "

$ws.Range("D31").Value = "'0.8628"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("D32").Value = "'5.093"
$ws.Range("E32").Value = "  -1.42%  "

$ws.Range("D33").Value = "'1.244"
$ws.Range("E33").Value = "  -4.38%  "

$ws.Range("D34").Value = "'3.057"
$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("D35").Value = "'0.05718"
$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("D36").Value = "'1.160"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").Value = "'0.9999"
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").Value = "'0.02040"
$ws.Range("E38").Value = "  -1.97%  "

$ws.Range("D39").Value = "'7.500"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").Value = "'0.5493"
$ws.Range("E40").Value = "  -2.72%  "

$ws.Range("D41").Value = "'2.969"
$ws.Range("E41").Value = "  +7.97%  "

$ws.Range("D42").Value = "'0.1751"
$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("D43").Value = "'9.372"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").Value = "'0.000002791"
$ws.Range("E44").Value = "  -7.52%  "

$ws.Range("D45").Value = "'2.174"
$ws.Range("E45").Value = "  +5.32%  "

$ws.Range("D46").Value = "'0.5178"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("D47").Value = "'11.25"
$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("D48").Value = "'0.06887"

$ws.Range("D49").Value = "'1.781"
$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").Value = "'110.39"
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("D51").Value = "'0.9997"
$ws.Range("E51").Value = "  -0.20%  "

